$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-09 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-10 Monday", 2) | Out-Null
$d.Content.Find.Execute("98÷8=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "54÷6=9, 0", 2) | Out-Null
$d.Content.Find.Execute("56÷8=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "47÷7=6, 5", 2) | Out-Null
$d.Content.Find.Execute("30÷3=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "94÷3=31, 1", 2) | Out-Null
$d.Content.Find.Execute("44÷2=22, 0", $true, $false, $false, $false, $false, $true, 1, $false, "94÷2=47, 0", 2) | Out-Null
$d.Content.Find.Execute("37÷7=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "82÷6=13, 4", 2) | Out-Null
$d.Content.Find.Execute("31÷6=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "39÷5=7, 4", 2) | Out-Null
$d.Content.Find.Execute("97÷3=32, 1", $true, $false, $false, $false, $false, $true, 1, $false, "56÷5=11, 1", 2) | Out-Null
$d.Content.Find.Execute("23÷3=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=6, 0", 2) | Out-Null
$d.Content.Find.Execute("87÷5=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=9, 3", 2) | Out-Null
$d.Content.Find.Execute("32÷3=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "74÷8=9, 2", 2) | Out-Null
$d.Content.Find.Execute("43÷5=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "45÷8=5, 5", 2) | Out-Null
$d.Content.Find.Execute("93÷3=31, 0", $true, $false, $false, $false, $false, $true, 1, $false, "68÷5=13, 3", 2) | Out-Null
$d.Content.Find.Execute("95÷2=47, 1", $true, $false, $false, $false, $false, $true, 1, $false, "59÷2=29, 1", 2) | Out-Null
$d.Content.Find.Execute("43÷8=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "43÷2=21, 1", 2) | Out-Null
$d.Content.Find.Execute("12÷8=1, 4", $true, $false, $false, $false, $false, $true, 1, $false, "58÷6=9, 4", 2) | Out-Null
$d.Content.Find.Execute("55÷7=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "88÷2=44, 0", 2) | Out-Null
$d.Content.Find.Execute("51÷5=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "62÷4=15, 2", 2) | Out-Null
$d.Content.Find.Execute("86÷4=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=12, 3", 2) | Out-Null
$d.Content.Find.Execute("77÷5=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "22÷4=5, 2", 2) | Out-Null
$d.Content.Find.Execute("77÷8=9, 5", $true, $false, $false, $false, $false, $true, 1, $false, "52÷2=26, 0", 2) | Out-Null
$d.Content.Find.Execute("77÷7=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "11÷5=2, 1", 2) | Out-Null
$d.Content.Find.Execute("86÷5=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "26÷3=8, 2", 2) | Out-Null
$d.Content.Find.Execute("84÷6=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "53÷3=17, 2", 2) | Out-Null
$d.Content.Find.Execute("21÷4=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "73÷7=10, 3", 2) | Out-Null
$d.Content.Find.Execute("64÷9=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "81÷8=10, 1", 2) | Out-Null
